$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FINAL DATA")

# Update error-bar half-widths in H6:I7 to be derived from the new bootstrapped
# confidence intervals (formulas instead of the old static values).
$ws.Range("H6").Formula = "=(8.91-8.1)/2"
$ws.Range("I6").Formula = "=(9.75-8.87)/2"
$ws.Range("H7").Formula = "=(9.5-8.63)/2"
$ws.Range("I7").Formula = "=(9.95-9.16)/2"

# New section further down the sheet documenting the bootstrapped CIs (pasted
# R console output) that back the new error bars above.
$ws.Range("B18").Value = '> quantile(PersonSOV.boot.mean$thetastar, c(0.025, 0.975))'
$ws.Range("B19").Value = "    2.5%    97.5% "
$ws.Range("B20").Value = "8.104762 8.914286 "
$ws.Range("B21").Value = '> PersonSVO.boot.mean = bootstrap(Scores[Scores$sentType=="AO" & Scores$sentOrder=="SVO",]$CorrectScore, 1000, mean)'
$ws.Range("B22").Value = '> quantile(PersonSVO.boot.mean$thetastar, c(0.025, 0.975))'
$ws.Range("B23").Value = "    2.5%    97.5% "
$ws.Range("B24").Value = "8.866667 9.752381 "
$ws.Range("B25").Value = '> ObjectSOV.boot.mean = bootstrap(Scores[Scores$sentType=="IO" & Scores$sentOrder=="SOV",]$CorrectScore, 1000, mean)'
$ws.Range("B26").Value = '> quantile(ObjectSOV.boot.mean$thetastar, c(0.025, 0.975))'
$ws.Range("B27").Value = "    2.5%    97.5% "
$ws.Range("B28").Value = "8.628571 9.505238 "
$ws.Range("B29").Value = '> ObjectSVO.boot.mean = bootstrap(Scores[Scores$sentType=="IO" & Scores$sentOrder=="SVO",]$CorrectScore, 1000, mean)'
$ws.Range("B30").Value = '> quantile(ObjectSVO.boot.mean$thetastar, c(0.025, 0.975))'
$ws.Range("B31").Value = "    2.5%    97.5% "
$ws.Range("B32").Value = "9.161905 9.952381 "

$ws.Range("B16").Value = "New error bars!  Bootstrapped Cis"

$ws.Range("E13").Select()
